$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CasesTab query (B2) to include the Cohort column
$casesQueryWithCohort = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['MGT01'] and demo.breed in ['Australian Cattle Dog','Mixed Breed']and diag.disease_term in ['Mammary Cancer'] and diag.primary_disease_site in ['Mammary Gland']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment` 
        coalesce(co.cohort_description, '') AS `Cohort`
'@
$ws.Range("B2").Value = $casesQueryWithCohort

# New row 5: StudyFilesTab
$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['MGT01'] and demo.breed in ['Australian Cattle Dog','Mixed Breed']and diag.disease_term in ['Mammary Cancer'] and diag.primary_disease_site in ['Mammary Gland']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(f.file_type, '') AS `File Type`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `File Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@
$studyStatQuery = @'
MATCH (s:study)
  MATCH (demo:demographic) 
  MATCH (diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['MGT01'] and demo.breed in ['Australian Cattle Dog','Mixed Breed']and diag.disease_term in ['Mammary Cancer'] and diag.primary_disease_site in ['Mammary Gland']
    
OPTIONAL MATCH (s)<-[:member_of]-(c:case)
OPTIONAL MATCH (c)<-[:of_case]-(samp:sample)<-[:of_sample]-(f:file)
RETURN 
	count(DISTINCT(f)) as number_of_files , 
	count(DISTINCT(samp)) as number_of_sample , 
	count(DISTINCT(c.case_id)) as number_of_cases , 
	count(DISTINCT(s.clinical_study_designation)) as number_of_study
'@

$ws.Range("A5").Value = "StudyFilesTab"
$ws.Range("B5").Value = $filesQuery
$ws.Range("C5").Value = $studyStatQuery
$ws.Range("D5").Value = "TC01_Canine_StudyMGT-Breed_Diagnosis_PrimDiseaseSite_Neo4jData.xlsx"
$ws.Range("E5").Value = "TC01_Canine_StudyMGT-Breed_Diagnosis_PrimDiseaseSite_WebData.xlsx"

# Apply the wrap-text style used by the other query cells (B:C) to the new row
$ws.Range("B5:C5").WrapText = $true

# Row heights
$ws.Rows.Item(2).RowHeight = 304.5
$ws.Rows.Item(5).RowHeight = 261

# Column widths (A widened slightly for "StudyFilesTab", E widened for longer file names wrap)
$ws.Columns.Item(1).ColumnWidth = 11.346354166666666
$ws.Columns.Item(5).ColumnWidth = 62.709635416666664

# Update the selected cell to C5
$ws.Range("C5").Select()
